$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.62"
$ws.Range("E2").Value = "'1.90%"

$ws.Range("D3").Value = "'41.79"
$ws.Range("E3").Value = "'3.46%"

$ws.Range("D4").Value = "'5.012"
$ws.Range("E4").Value = "'-0.09%"

$ws.Range("D5").Value = "'0.07523"
$ws.Range("E5").Value = "'3.13%"

$ws.Range("D6").Value = "'1.595"
$ws.Range("E6").Value = "'4.26%"

$ws.Range("D7").Value = "'0.9153"
$ws.Range("E7").Value = "'-1.18%"

$ws.Range("D9").Value = "'0.1179"
$ws.Range("E9").Value = "'2.01%"

$ws.Range("D10").Value = "'0.1823"
$ws.Range("E10").Value = "'3.45%"

$ws.Range("D11").Value = "'0.08879"
$ws.Range("E11").Value = "'1.09%"

$ws.Range("D12").Value = "'0.04109"
$ws.Range("E12").Value = "'-5.68%"

$ws.Range("D13").Value = "'0.1050"
$ws.Range("E13").Value = "'-0.24%"

$ws.Range("D14").Value = "'0.001290"
$ws.Range("E14").Value = "'1.84%"

$ws.Range("E15").Value = "'1.18%"

$ws.Range("D16").Value = "'3.345"
$ws.Range("E16").Value = "'0.17%"

$ws.Range("D17").Value = "'4.362"
$ws.Range("E17").Value = "'1.80%"

$ws.Range("D18").Value = "'0.3329"
$ws.Range("E18").Value = "'1.50%"

$ws.Range("D19").Value = "'8.296"
$ws.Range("E19").Value = "'4.01%"

$ws.Range("E20").Value = "'-2.82%"

$ws.Range("E21").Value = "'11.80%"

$ws.Range("D22").Value = "'0.04099"
$ws.Range("E22").Value = "'4.38%"

$ws.Range("E23").Value = "'0.32%"

$ws.Range("D24").Value = "'0.003887"
$ws.Range("E24").Value = "'2.82%"

$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'8.36%"

$ws.Range("D38").Value = "'0.02390"
$ws.Range("E38").Value = "'4.09%"

$ws.Range("D39").Value = "'0.05214"
$ws.Range("E39").Value = "'3.42%"

$ws.Range("D40").Value = "'0.006991"
$ws.Range("E40").Value = "'18.09%"

$ws.Range("D41").Value = "'0.007775"
$ws.Range("E41").Value = "'-0.91%"

$ws.Range("D42").Value = "'0.1323"
$ws.Range("E42").Value = "'3.07%"

$ws.Range("D43").Value = "'0.007418"
$ws.Range("E43").Value = "'0.75%"

$ws.Range("D44").Value = "'0.007698"
$ws.Range("E44").Value = "'6.40%"

$ws.Range("D45").Value = "'0.3238"
$ws.Range("E45").Value = "'1.54%"

$ws.Range("D46").Value = "'0.00006593"
$ws.Range("E46").Value = "'6.77%"

$ws.Range("D48").Value = "'0.04533"
$ws.Range("E48").Value = "'16.54%"

$ws.Range("D49").Value = "'0.004206"
$ws.Range("E49").Value = "'0.11%"

$ws.Range("D50").Value = "'0.00002103"

$ws.Range("D51").Value = "'0.0002003"
